$d = $word.ActiveDocument

$replacements = @(
    @("2026-01-04 Sunday", "2026-01-05 Monday"),
    @("822÷4=205, 2", "310÷8=38, 6"),
    @("626÷8=78, 2", "310÷2=155, 0"),
    @("795÷5=159, 0", "914÷2=457, 0"),
    @("284÷3=94, 2", "820÷2=410, 0"),
    @("997÷7=142, 3", "642÷9=71, 3"),
    @("211÷3=70, 1", "514÷9=57, 1"),
    @("971÷5=194, 1", "570÷2=285, 0"),
    @("854÷8=106, 6", "999÷8=124, 7"),
    @("910÷5=182, 0", "584÷3=194, 2"),
    @("270÷8=33, 6", "577÷2=288, 1"),
    @("207÷5=41, 2", "107÷3=35, 2"),
    @("762÷3=254, 0", "397÷7=56, 5"),
    @("166÷2=83, 0", "606÷7=86, 4"),
    @("894÷4=223, 2", "247÷9=27, 4"),
    @("696÷2=348, 0", "484÷7=69, 1"),
    @("247÷5=49, 2", "603÷2=301, 1"),
    @("110÷6=18, 2", "648÷5=129, 3"),
    @("281÷4=70, 1", "642÷7=91, 5"),
    @("595÷4=148, 3", "166÷9=18, 4"),
    @("129÷4=32, 1", "345÷9=38, 3"),
    @("301÷3=100, 1", "553÷7=79, 0"),
    @("716÷2=358, 0", "981÷3=327, 0"),
    @("426÷9=47, 3", "445÷3=148, 1"),
    @("182÷4=45, 2", "175÷7=25, 0"),
    @("296÷8=37, 0", "388÷3=129, 1")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Write-Output "Done: applied $($replacements.Count) replacements"
